$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$old = "dnasr281@gmail.com, System"
$new = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $old) {
        $cell.Value = $new
    }
}
